$d = $word.ActiveDocument

# --- Update the date/day heading (unique text in the document) ---
$d.Content.Find.Execute("2024-02-23 Friday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2024-02-24 Saturday", 2)

# --- Update the division problems in the table ---
# The table has 5 "data" rows (1, 5, 9, 13, 17) each with 5 columns of problems;
# the rows in between are blank spacer rows. We address each problem cell
# directly by (row, col) and overwrite its Range.Text, which avoids any
# ambiguity from duplicate problem text elsewhere in the table (e.g. "68÷7="
# appears twice with two different replacements).
$t = $d.Tables.Item(1)

$edits = @(
    @{ Row = 1;  Col = 1; New = "53÷8=" },
    @{ Row = 1;  Col = 2; New = "94÷7=" },
    @{ Row = 1;  Col = 3; New = "99÷8=" },
    @{ Row = 1;  Col = 4; New = "83÷8=" },
    @{ Row = 1;  Col = 5; New = "21÷3=" },

    @{ Row = 5;  Col = 1; New = "83÷7=" },
    @{ Row = 5;  Col = 2; New = "42÷7=" },
    @{ Row = 5;  Col = 3; New = "95÷8=" },
    @{ Row = 5;  Col = 4; New = "36÷2=" },
    @{ Row = 5;  Col = 5; New = "62÷3=" },

    @{ Row = 9;  Col = 1; New = "93÷3=" },
    @{ Row = 9;  Col = 2; New = "58÷4=" },
    @{ Row = 9;  Col = 3; New = "80÷3=" },
    @{ Row = 9;  Col = 4; New = "77÷2=" },
    @{ Row = 9;  Col = 5; New = "25÷5=" },

    @{ Row = 13; Col = 1; New = "61÷4=" },
    @{ Row = 13; Col = 2; New = "40÷3=" },
    @{ Row = 13; Col = 3; New = "42÷8=" },
    @{ Row = 13; Col = 4; New = "61÷3=" },
    @{ Row = 13; Col = 5; New = "57÷8=" },

    @{ Row = 17; Col = 1; New = "48÷5=" },
    @{ Row = 17; Col = 2; New = "22÷7=" },
    @{ Row = 17; Col = 3; New = "60÷7=" },
    @{ Row = 17; Col = 4; New = "36÷3=" },
    @{ Row = 17; Col = 5; New = "78÷3=" }
)

foreach ($e in $edits) {
    $t.Cell($e.Row, $e.Col).Range.Text = $e.New
}
